$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "21×40="
$t.Cell(1,2).Range.Text = "91×64="
$t.Cell(1,3).Range.Text = "84×58="
$t.Cell(1,4).Range.Text = "84×75="
$t.Cell(1,5).Range.Text = "19×52="
$t.Cell(5,1).Range.Text = "53×49="
$t.Cell(5,2).Range.Text = "82×90="
$t.Cell(5,3).Range.Text = "90×44="
$t.Cell(5,4).Range.Text = "76×91="
$t.Cell(5,5).Range.Text = "49×21="
$t.Cell(10,1).Range.Text = "88×23="
$t.Cell(10,2).Range.Text = "85×96="
$t.Cell(10,3).Range.Text = "93×27="
$t.Cell(10,4).Range.Text = "75×62="
$t.Cell(10,5).Range.Text = "88×64="
$t.Cell(15,1).Range.Text = "58×84="
$t.Cell(15,2).Range.Text = "95×60="
$t.Cell(15,3).Range.Text = "67×46="
$t.Cell(15,4).Range.Text = "79×72="
$t.Cell(15,5).Range.Text = "50×61="
$t.Cell(20,1).Range.Text = "71×72="
$t.Cell(20,2).Range.Text = "75×87="
$t.Cell(20,3).Range.Text = "17×13="
$t.Cell(20,4).Range.Text = "75×57="
$t.Cell(20,5).Range.Text = "67×81="
